# Apply changes described by the commit:
# "changed stimulus duration to 5 s" -> update stimulus_duration column (E) for rows 2-4
# Selection moved to F4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update stimulus_duration (column E) values for rows 2, 3, 4 from 2 -> 5
$ws.Range("E2").Value = 5
$ws.Range("E3").Value = 5
$ws.Range("E4").Value = 5

# Update the active selection to F4 as reflected in the saved sheet view
$ws.Range("F4").Select()
